$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.46425193330187
$ws.Range("C2").Value = 0.0000000000000001110223024625157
$ws.Range("D2").Value = 0.01533359077428826
$ws.Range("E2").Value = 0.9291784833936778
$ws.Range("F2").Value = 0.8633726540017753

$ws.Range("B3").Value = 16.4839847028106
$ws.Range("C3").Value = 0.0000000000000001110223024625157
$ws.Range("D3").Value = 0.01530496395580449
$ws.Range("E3").Value = 0.9274437674895769
$ws.Range("F3").Value = 0.8601519418552606

$ws.Range("B4").Value = 17.34811180985909
$ws.Range("C4").Value = 0.0000000000000001110223024625157
$ws.Range("D4").Value = 0.01299929766748106
$ws.Range("E4").Value = 0.7877259716691299
$ws.Range("F4").Value = 0.6205122064420749

$ws.Range("B5").Value = 17.68010677920795
$ws.Range("C5").Value = 0.0000000000000001110223024625157
$ws.Range("D5").Value = 0.01347432114524307
$ws.Range("E5").Value = 0.8165112445475096
$ws.Range("F5").Value = 0.6666906124725229

$ws.Range("B6").Value = 17.34630636150452
$ws.Range("C6").Value = 0.0000000000000001110223024625157
$ws.Range("D6").Value = 0.01298025924589215
$ws.Range("E6").Value = 0.7865722894065345
$ws.Range("F6").Value = 0.6186959664622371
